$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 234. This shifts the existing row 234 (and everything
# below it) down by one, matching the diff: old row 234 -> new row 235,
# old row 235 -> new row 236, ..., old row 281 -> new row 282.
$ws.Rows.Item(234).Insert()

# Fill the newly inserted row 234 with the new record (a "Larry Ann" lot
# reported for this market/product on a later date).
$ws.Range("A234").Value2 = 4
$ws.Range("B234").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C234").Value2 = "Los Lagos"
$ws.Range("D234").Value2 = 44964
$ws.Range("E234").Value2 = 10
$ws.Range("F234").Value2 = "Fruta"
$ws.Range("G234").Value2 = 100103
$ws.Range("H234").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I234").Value2 = 100103002
$ws.Range("J234").Value2 = "Ciruela"
$ws.Range("K234").Value2 = "Larry Ann"
$ws.Range("L234").Value2 = "Primera"
$ws.Range("M234").Value2 = 600
$ws.Range("N234").Value2 = 17000
$ws.Range("O234").Value2 = 18000
$ws.Range("P234").Value2 = 17500
$ws.Range("Q234").Value2 = "$/caja 14 kilos granel"
$ws.Range("R234").Value2 = "Región de O'Higgins"
$ws.Range("S234").Value2 = 1250
$ws.Range("T234").Value2 = 14
